# Apply updated crypto price/volume data (GitHub Actions scrape refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.896.04'
$ws.Range('E2').Value = '  -0.41%  '
$ws.Range('D3').Value = '2.556.39'
$ws.Range('E3').Value = '  +0.29%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.13%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '305.05'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.21%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '99.00'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +6.60%  '
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('E9').Value = '  -0.31%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '37.12'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +3.15%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0809'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.12%  '
$ws.Range('B12').Value = 'Polkadot'
$ws.Range('C12').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.77'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.05%  '
$ws.Range('B13').Value = 'TRON'
$ws.Range('C13').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.116'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +7.09%  '
$ws.Range('D14').Value = '2.519.43'
$ws.Range('E14').Value = '  -2.11%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.06'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +5.95%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.884'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.36%  '
$ws.Range('D17').Value = '42.939.69'
$ws.Range('E17').Value = '  -0.44%  '
$ws.Range('E18').Value = '  +7.23%  '
$ws.Range('E19').Value = '  +0.65%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.65'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.43%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '71.75'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.11%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '255.63'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.20%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.99'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.27%  '
$ws.Range('E24').Value = '  -1.96%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '27.87'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -5.61%  '
$ws.Range('E26').Value = '  -0.03%  '
$ws.Range('E27').Value = '  +0.62%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '38.32'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +3.66%  '
$ws.Range('E30').Value = '  +0.69%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '158.44'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.73%  '
$ws.Range('E32').Value = '  +0.08%  '
$ws.Range('E33').Value = '  +0.22%  '
$ws.Range('E34').Value = '  +1.23%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.33'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -2.15%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '26.65'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +13.12%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '18.70'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +13.20%  '
$ws.Range('E38').Value = '  -1.21%  '
$ws.Range('E39').Value = '  -0.32%  '
$ws.Range('B40').Value = 'NEARProtocol'
$ws.Range('C40').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.50'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.53%  '
$ws.Range('B41').Value = 'ApeXProtocol'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.12'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +33.67%  '
$ws.Range('E42').Value = '  +0.07%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0307'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.16%  '
$ws.Range('D44').Value = '2.085.01'
$ws.Range('E44').Value = '  +0.15%  '
$ws.Range('E45').Value = '  +0.00%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '86.99'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.24%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.06'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +3.39%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '75.61'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +8.61%  '
$ws.Range('D49').Value = '2.805.36'
$ws.Range('E49').Value = '  +0.26%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '104.12'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.27%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.191'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.44%  '

Write-Output "Updated crypto price/volume data"
